# CS371ReportDocument.docx edit
#
# The commit's only substantive change to visible document text is in the
# "Implementation:" line: the sentence describing where the diagram is
# submitted was expanded/reworded. (All of the other hunks in the source
# diff only add <w:proofErr> spell/grammar-check bookkeeping tags and shift
# a <w:lastRenderedPageBreak/> pagination hint around paragraphs whose
# run text is otherwise unchanged - those are artifacts Word's own
# background proofing/pagination pass stamps into the XML and are not
# reachable through the Word object model, so there is nothing for this
# script to "do" for those hunks beyond leaving the text exactly as-is.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Diagram is in submitted files on Canvas",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Diagram is in the files that are zipped together and submitted on Canvas",
    2
)
